$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "estado de cuenta" detail rows (16-24) with the new data set.
# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico

$ws.Range("C16").Value = "1143382942"
$ws.Range("D16").Value = "LUIS CARLOS MEDRANO PERNETH"
$ws.Range("E16").Value = "1904"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

$ws.Range("C17").Value = "1143391777"
$ws.Range("D17").Value = "MARIA JOSE BERMEJO SALGUEDO"
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 53653
$ws.Range("G17").Value = 1392310

$ws.Range("C18").Value = "1001967915"
$ws.Range("D18").Value = "CAMILO ANDRES ORTEGA FUENTES"
$ws.Range("E18").Value = "1907"
$ws.Range("F18").Value = 53653
$ws.Range("G18").Value = 1341340

$ws.Range("C19").Value = "1143391777"
$ws.Range("D19").Value = "MARIA JOSE BERMEJO SALGUEDO"
$ws.Range("E19").Value = "1908"
$ws.Range("F19").Value = 53653
$ws.Range("G19").Value = 1392310

$ws.Range("C20").Value = "1001967915"
$ws.Range("D20").Value = "CAMILO ANDRES ORTEGA FUENTES"
$ws.Range("E20").Value = "1908"
$ws.Range("F20").Value = 53653
$ws.Range("G20").Value = 1341340

$ws.Range("C21").Value = "1143391777"
$ws.Range("D21").Value = "MARIA JOSE BERMEJO SALGUEDO"
$ws.Range("E21").Value = "1909"
$ws.Range("F21").Value = 53653
$ws.Range("G21").Value = 1392310

$ws.Range("C22").Value = "1143391777"
$ws.Range("D22").Value = "MARIA JOSE BERMEJO SALGUEDO"
$ws.Range("E22").Value = "1910"
$ws.Range("F22").Value = 53653
$ws.Range("G22").Value = 1392310

$ws.Range("C23").Value = "1143391777"
$ws.Range("D23").Value = "MARIA JOSE BERMEJO SALGUEDO"
$ws.Range("E23").Value = "1911"
$ws.Range("F23").Value = 53653
$ws.Range("G23").Value = 1392310

$ws.Range("C24").Value = "1047479968"
$ws.Range("D24").Value = "ROSANGELA ORDOSGOITIA MOHADIE"
$ws.Range("E24").Value = "2001"
$ws.Range("F24").Value = 53653
$ws.Range("G24").Value = 1341340

# Nudge the logo image to the left to match the updated layout, keeping its
# exact original size.
$shp = $ws.Shapes.Item(1)
$shp.Width = 76.81889763779527
$shp.Height = 48.188976377952756
$shp.Left = 59.09055118110236
$shp.Top = 19.405511811023622

# Re-tighten the detail columns now that the refreshed data no longer needs
# as much horizontal room (mirrors the author's column narrowing pass).
$ws.Columns.Item(2).ColumnWidth = 15.916916666666667
$ws.Columns.Item(3).ColumnWidth = 9.917072916666665
$ws.Columns.Item(4).ColumnWidth = 30.583635416666667
$ws.Columns.Item(5).ColumnWidth = 11.751229166666667
$ws.Columns.Item(6).ColumnWidth = 8.583791666666666
$ws.Columns.Item(7).ColumnWidth = 12.583791666666666
$ws.Columns.Item(8).ColumnWidth = 16.91691666666667
$ws.Columns.Item(9).ColumnWidth = 15.917072916666665
$ws.Columns.Item(10).ColumnWidth = 13.250354166666666
